$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the four new task rows' column headers (C2:F2) for Wiki update / Log4j /
# Junit tests / Hibernate, matching the new shared-string entries.
$ws.Range("C2").Value = "Wiki update"
$ws.Range("D2").Value = "Log4j"
$ws.Range("E2").Value = "Junit tests"
$ws.Range("F2").Value = "Hibernate"

# Auto-fit the new columns to their content, same as columns A and B already are.
$ws.Columns.Item(3).ColumnWidth = 10.833333333333332
$ws.Columns.Item(4).ColumnWidth = 4.666666666666667
$ws.Columns.Item(5).ColumnWidth = 9.166666666666666
$ws.Columns.Item(6).ColumnWidth = 9.0

# Move the active selection the way the author left it.
$ws.Range("C3").Select()
